# Applies the cryptos list refresh described in the commit diff.
# For each touched cell we set the literal text value. Columns D may contain
# strings that Excel would otherwise auto-parse as numbers (e.g. "584.81"),
# so for those we use a leading apostrophe to force text, then reset the
# cell style back to Normal so no stray number-format style is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "63.222.43"
$ws.Range("E2").Value = "  +6.07%  "
# Row 3
$ws.Range("D3").Value = "3.109.03"
$ws.Range("E3").Value = "  +3.95%  "
# Row 4
$ws.Range("E4").Value = "  +0.03%  "
# Row 5
$ws.Range("D5").Value = "'584.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.12%  "
# Row 6
$ws.Range("D6").Value = "'144.87"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.18%  "
# Row 7
$ws.Range("E7").Value = "  +0.00%  "
# Row 8
$ws.Range("D8").Value = "3.102.33"
$ws.Range("E8").Value = "  +4.00%  "
# Row 9
$ws.Range("D9").Value = "'0.531"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +2.11%  "
# Row 10
$ws.Range("D10").Value = "'0.152"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +13.76%  "
# Row 11
$ws.Range("E11").Value = "  +7.16%  "
# Row 12
$ws.Range("D12").Value = "'0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +3.58%  "
# Row 13
$ws.Range("E13").Value = "  +7.43%  "
# Row 14
$ws.Range("D14").Value = "'35.41"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.68%  "
# Row 15
$ws.Range("E15").Value = "  +0.66%  "
# Row 16
$ws.Range("D16").Value = "3.623.14"
$ws.Range("E16").Value = "  +3.99%  "
# Row 17
$ws.Range("E17").Value = "  +1.00%  "
# Row 18
$ws.Range("D18").Value = "63.128.04"
$ws.Range("E18").Value = "  +5.96%  "
# Row 19
$ws.Range("D19").Value = "3.108.99"
$ws.Range("E19").Value = "  +4.11%  "
# Row 20
$ws.Range("D20").Value = "'466.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +6.40%  "
# Row 21
$ws.Range("E21").Value = "  +3.57%  "
# Row 22
$ws.Range("E22").Value = "  +0.23%  "
# Row 23
$ws.Range("E23").Value = "  +6.84%  "
# Row 24
$ws.Range("D24").Value = "'13.30"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.29%  "
# Row 25
$ws.Range("D25").Value = "'81.98"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.34%  "
# Row 26
$ws.Range("E26").Value = "  -0.06%  "
# Row 27
$ws.Range("D27").Value = "'8.55"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +10.22%  "
# Row 28
$ws.Range("E28").Value = "  -0.43%  "
# Row 29
$ws.Range("B29").Value = "PancakeSwap"
$ws.Range("C29").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D29").Value = "'2.67"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +4.37%  "
# Row 30
$ws.Range("B30").Value = "FirstDigitalUSD"
$ws.Range("C30").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D30").Value = "'1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.07%  "
# Row 31
$ws.Range("E31").Value = "  +9.93%  "
# Row 32
$ws.Range("D32").Value = "'26.94"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.10%  "
# Row 33
$ws.Range("E33").Value = "  +4.03%  "
# Row 34
$ws.Range("D34").Value = "0.0₃0866"
$ws.Range("E34").Value = "  +11.12%  "
# Row 35
$ws.Range("D35").Value = "'2.41"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +15.87%  "
# Row 36
$ws.Range("E36").Value = "  +6.47%  "
# Row 37
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +20.42%  "
# Row 38
$ws.Range("B38").Value = "Filecoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D38").Value = "'6.08"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.63%  "
# Row 39
$ws.Range("D39").Value = "'50.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.90%  "
# Row 40
$ws.Range("D40").Value = "'435.96"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.10%  "
# Row 41
$ws.Range("E41").Value = "  +0.30%  "
# Row 42
$ws.Range("D42").Value = "2.912.91"
$ws.Range("E42").Value = "  +6.42%  "
# Row 43
$ws.Range("E43").Value = "  +4.54%  "
# Row 44
$ws.Range("E44").Value = "  +11.19%  "
# Row 45
$ws.Range("E45").Value = "  +4.46%  "
# Row 46
$ws.Range("E46").Value = "  +7.23%  "
# Row 48
$ws.Range("D48").Value = "'34.90"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.15%  "
# Row 49
$ws.Range("D49").Value = "'122.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.59%  "
# Row 50
$ws.Range("E50").Value = "  +0.96%  "
# Row 51
$ws.Range("D51").Value = "'24.61"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.30%  "
